# "android to android pending"
# Update a set of "Yes/No" <-> "Yes" cells on the Test List sheet, and move
# the sheet's current view/selection down to where the Android rows are.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet2")

# Cells that flip from "Yes/No" to "Yes"
$toYes = @("M6","M8","K10","L10","M10","E16","E18","E20","C22","D22","E22","E28","E30","E32","C34","D34","E34","E40","E42","E44","C46","D46","E46")
foreach ($addr in $toYes) {
    $ws.Range($addr).Value = "Yes"
}

# Cell that flips from "Yes" to "Yes/No" (Android -> Android, Discovery block)
$ws.Range("B10").Value = "Yes/No"

# Update the view: scroll the window down and move the selection to F47
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("F47").Select()
